$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "FT2325009V3R2J01"
